# Update the "Attributes" column (column D) values on the active worksheet
# to reflect the regenerated skeleton-code attribute sets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "{'doctor_licence_number', 'Doctor_Name', 'hospital_department'},{'List_plan', 'Name_of_clinics', 'Preventive_care'}"
$ws.Range("D3").Value  = "{'List_plan', 'Name_of_clinics', 'Preventive_care'}{'Patient_age', 'Schedule_time', 'Name_of_clinics', 'Customer_phone'},{'Preventive_care', 'Address_Of_clinics', 'Policy_number'}"
$ws.Range("D4").Value  = "{'Acc_type', 'Address_Of_clinics', 'Hospital_Address', 'Preventive_care'}"
$ws.Range("D5").Value  = "{'doctor_licence_number', 'Doctor_Name', 'hospital_department'}{'doctor_licence_number', 'Coverage_policy', 'Doctor_Name', 'Hourly_charge_doctor', 'hospital_department'},{'Customer_phone', 'Doctor_Name', 'Hospital_Address', 'Patient_prior_condition', 'Patient_age'}"
$ws.Range("D6").Value  = "{'List_plan', 'Name_of_clinics', 'Preventive_care'}"
$ws.Range("D7").Value  = "{'X,Y_Coordinates'}{'Customer_phone', 'Doctor_Name', 'Hospital_Address', 'Patient_prior_condition', 'Patient_age'}"
$ws.Range("D8").Value  = "{'Acc_type', 'Schedule_time', 'Doctor_Name', 'Hospital_Address', 'Discharge_amount'}"
$ws.Range("D9").Value  = "{'Schedule_time', 'Doctor_available_time'}"
$ws.Range("D10").Value = "{'Schedule_time', 'Doctor_available_time'}{'doctor_licence_number', 'Doctor_Name', 'hospital_department'},{'doctor_licence_number', 'Coverage_policy', 'Doctor_Name', 'Hourly_charge_doctor', 'hospital_department'}{'Schedule_time', 'Doctor_available_time'}"
